# Weekly update: insert a new daily price record for "Uva" (Agrícola del
# Norte S.A. de Arica) as the new row 42, pushing the existing rows
# 42:70 down to 43:71 (dimension grows from A1:T70 to A1:T71).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 42 - this shifts rows 42:70
# down to 43:71 and extends the sheet dimension automatically.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new record.
$ws.Range("A42").Value = 1
$ws.Range("B42").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C42").Value = "Arica y Parinacota"
$ws.Range("D42").Value = 44572
$ws.Range("E42").Value = 15
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100109
$ws.Range("H42").Value = "Uva"
$ws.Range("I42").Value = 100109001
$ws.Range("J42").Value = "Uva"
$ws.Range("K42").Value = "Superior Seedless"
$ws.Range("L42").Value = "Segunda"
$ws.Range("M42").Value = 250
$ws.Range("N42").Value = 13000
$ws.Range("O42").Value = 14000
$ws.Range("P42").Value = 13500
$ws.Range("Q42").Value = "$/caja 10 kilos"
$ws.Range("R42").Value = "Región de O'Higgins"
$ws.Range("S42").Value = 1350
$ws.Range("T42").Value = 10
